$d = $word.ActiveDocument

# 1) Drop the "_GoBack" bookmark from the title paragraph ("Informe Técnico
#    de Actividades"). It will be re-created further below, at the end of
#    the document, where the new content is authored.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Turn the trailing empty paragraph into a new "Resultado" heading
#    (style "Ttulo1" — same style used by the document title).
$last = $d.Paragraphs.Last
$last.Style = "Ttulo1"
# A trailing sentinel character is used so a genuinely collapsed (empty)
# range can be captured for the bookmark below; it is stripped afterwards.
$last.Range.Text = "ResultadoX"

# 3) Re-create "_GoBack" positioned right after the new run's text (i.e.
#    collapsed at the end of "Resultado", matching how Word leaves it after
#    the last edit). We bookmark the sentinel character first (so the range
#    is non-degenerate) and then delete just that character — the
#    bookmark tags stay behind, now collapsed immediately after "Resultado".
$para = $d.Paragraphs.Last
$paraRange = $para.Range
$sentinelEnd = $paraRange.End - 1
$sentinelRange = $d.Range($sentinelEnd - 1, $sentinelEnd)
$d.Bookmarks.Add("_GoBack", $sentinelRange)

$sentinelRange = $d.Range($sentinelEnd - 1, $sentinelEnd)
$sentinelRange.Text = ""
